$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# All "Dimensions" cells in column F (rows 3-10 and 15-22 previously said
# "96x96"; rows 11-14 and 23-26 previously said "128x128"). The "128x128"
# shared string is being removed entirely, leaving "96x96" as the only
# dimension value, so every one of these cells now reads "96x96".
$ws.Range("F11:F14").Value = "96x96"
$ws.Range("F23:F26").Value = "96x96"

# Move the active selection from G7 to H14.
$ws.Range("H14").Select()
